$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension / add header for new column F by copying format from E1 (bordered, bold, centered header style)
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# time_taken values for rows 2-95 (column F), as recorded at panel export time
$timeTaken = @(
    "2021-10-05 13:40:51.376312",
    "2021-10-05 13:40:51.376323",
    "2021-10-05 13:40:51.376326",
    "2021-10-05 13:40:51.376329",
    "2021-10-05 13:40:51.376332",
    "2021-10-05 13:40:51.376334",
    "2021-10-05 13:40:51.376337",
    "2021-10-05 13:40:51.376339",
    "2021-10-05 13:40:51.376342",
    "2021-10-05 13:40:51.376345",
    "2021-10-05 13:40:51.376347",
    "2021-10-05 13:40:51.376349",
    "2021-10-05 13:40:51.376352",
    "2021-10-05 13:40:51.376354",
    "2021-10-05 13:40:51.376357",
    "2021-10-05 13:40:51.376360",
    "2021-10-05 13:40:51.376362",
    "2021-10-05 13:40:51.376365",
    "2021-10-05 13:40:51.376367",
    "2021-10-05 13:40:51.376370",
    "2021-10-05 13:40:51.376372",
    "2021-10-05 13:40:51.376374",
    "2021-10-05 13:40:51.376377",
    "2021-10-05 13:40:51.376379",
    "2021-10-05 13:40:51.376382",
    "2021-10-05 13:40:51.376385",
    "2021-10-05 13:40:51.376387",
    "2021-10-05 13:40:51.376390",
    "2021-10-05 13:40:51.376392",
    "2021-10-05 13:40:51.376395",
    "2021-10-05 13:40:51.376397",
    "2021-10-05 13:40:51.376400",
    "2021-10-05 13:40:51.376402",
    "2021-10-05 13:40:51.376405",
    "2021-10-05 13:40:51.376408",
    "2021-10-05 13:40:51.376410",
    "2021-10-05 13:40:51.376413",
    "2021-10-05 13:40:51.376415",
    "2021-10-05 13:40:51.376417",
    "2021-10-05 13:40:51.376420",
    "2021-10-05 13:40:51.376423",
    "2021-10-05 13:40:51.376425",
    "2021-10-05 13:40:51.376428",
    "2021-10-05 13:40:51.376430",
    "2021-10-05 13:40:51.376433",
    "2021-10-05 13:40:51.376435",
    "2021-10-05 13:40:51.376437",
    "2021-10-05 13:40:51.376440",
    "2021-10-05 13:40:51.376442",
    "2021-10-05 13:40:51.376445",
    "2021-10-05 13:40:51.376447",
    "2021-10-05 13:40:51.376450",
    "2021-10-05 13:40:51.376453",
    "2021-10-05 13:40:51.376455",
    "2021-10-05 13:40:51.376458",
    "2021-10-05 13:40:51.376460",
    "2021-10-05 13:40:51.376462",
    "2021-10-05 13:40:51.376465",
    "2021-10-05 13:40:51.376467",
    "2021-10-05 13:40:51.376470",
    "2021-10-05 13:40:51.376472",
    "2021-10-05 13:40:51.376475",
    "2021-10-05 13:40:51.376477",
    "2021-10-05 13:40:51.376480",
    "2021-10-05 13:40:51.376483",
    "2021-10-05 13:40:51.376486",
    "2021-10-05 13:40:51.376488",
    "2021-10-05 13:40:51.376491",
    "2021-10-05 13:40:51.376493",
    "2021-10-05 13:40:51.376496",
    "2021-10-05 13:40:51.376498",
    "2021-10-05 13:40:51.376501",
    "2021-10-05 13:40:51.376503",
    "2021-10-05 13:40:51.376506",
    "2021-10-05 13:40:51.376508",
    "2021-10-05 13:40:51.376511",
    "2021-10-05 13:40:51.376515",
    "2021-10-05 13:40:51.376518",
    "2021-10-05 13:40:51.376521",
    "2021-10-05 13:40:51.376524",
    "2021-10-05 13:40:51.376526",
    "2021-10-05 13:40:51.376528",
    "2021-10-05 13:40:51.376531",
    "2021-10-05 13:40:51.376533",
    "2021-10-05 13:40:51.376536",
    "2021-10-05 13:40:51.376538",
    "2021-10-05 13:40:51.376541",
    "2021-10-05 13:40:51.376544",
    "2021-10-05 13:40:51.376546",
    "2021-10-05 13:40:51.376549",
    "2021-10-05 13:40:51.376551",
    "2021-10-05 13:40:51.376554",
    "2021-10-05 13:40:51.376557",
    "2021-10-05 13:40:51.376560"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $timeTaken[$i]
}
